$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").Value = "67.447.97"
$ws.Range("E2").Value = "  +2.35%  "

# Row 3: update D3, E3
$ws.Range("D3").Value = "3.371.46"
$ws.Range("E3").Value = "  +1.49%  "

# Row 4: update E4
$ws.Range("E4").Value = "  -0.28%  "

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.29"
$ws.Range("E5").Value = "  +6.48%  "

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.35"
$ws.Range("E6").Value = "  -0.51%  "

# Row 7: update D7, E7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8: update D8, E8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +2.00%  "

# Row 9: update E9
$ws.Range("E9").Value = "  +1.81%  "

# Row 10: update D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.588"
$ws.Range("E10").Value = "  +1.32%  "

# Row 11: update D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "47.45"
$ws.Range("E11").Value = "  +2.35%  "

# Row 12: update D12, E12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000274"
$ws.Range("E12").Value = "  +2.41%  "

# Row 13: update D13, E13
$ws.Range("D13").Value = "3.919.67"
$ws.Range("E13").Value = "  +1.17%  "

# Row 14: update B14, C14, D14, E14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.61"
$ws.Range("E14").Value = "  +0.53%  "

# Row 15: update B15, C15, D15, E15
$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "634.24"
$ws.Range("E15").Value = "  +7.16%  "

# Row 16: update D16, E16
$ws.Range("D16").Value = "67.601.20"
$ws.Range("E16").Value = "  +2.12%  "

# Row 17: update B17, C17, D17, E17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.119"
$ws.Range("E17").Value = "  +1.33%  "

# Row 18: update B18, C18, D18, E18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.375.16"
$ws.Range("E18").Value = "  +1.20%  "

# Row 19: update D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.00"
$ws.Range("E19").Value = "  +0.78%  "

# Row 20: update D20, E20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.16"
$ws.Range("E20").Value = "  +1.80%  "

# Row 21: update D21, E21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.910"
$ws.Range("E21").Value = "  +1.55%  "

# Row 22: update D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.99"
$ws.Range("E22").Value = "  -3.06%  "

# Row 23: update E23
$ws.Range("E23").Value = "  +2.14%  "

# Row 24: update D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "99.29"
$ws.Range("E24").Value = "  +0.47%  "

# Row 25: update D25, E25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.00"
$ws.Range("E25").Value = "  +1.45%  "

# Row 26: update E26
$ws.Range("E26").Value = "  +4.58%  "

# Row 27: update D27, E27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.69"
$ws.Range("E27").Value = "  +2.36%  "

# Row 28: update D28, E28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "32.58"
$ws.Range("E28").Value = "  +7.44%  "

# Row 29: update D29, E29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.68"
$ws.Range("E29").Value = "  +1.41%  "

# Row 30: update D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.95"
$ws.Range("E30").Value = "  +4.20%  "

# Row 31: update D31, E31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "609.17"
$ws.Range("E31").Value = "  +5.44%  "

# Row 32: update D32, E32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.80"
$ws.Range("E32").Value = "  -2.77%  "

# Row 33: update D33, E33
$ws.Range("D33").Value = "3.999.98"
$ws.Range("E33").Value = "  +7.98%  "

# Row 34: update D34, E34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.09"
$ws.Range("E34").Value = "  +1.47%  "

# Row 35: update E35
$ws.Range("E35").Value = "  +1.92%  "

# Row 36: update D36, E36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.13%  "

# Row 37: update D37, E37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.00"
$ws.Range("E37").Value = "  -0.11%  "

# Row 38: update D38, E38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.83"
$ws.Range("E38").Value = "  +7.11%  "

# Row 39: update D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.132"
$ws.Range("E39").Value = "  +4.17%  "

# Row 40: update D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.70"
$ws.Range("E40").Value = "  +0.63%  "

# Row 41: update D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.24"
$ws.Range("E41").Value = "  +2.04%  "

# Row 42: update D42, E42
$ws.Range("D42").Value = "0.0₃0702"
$ws.Range("E42").Value = "  +0.49%  "

# Row 43: update E43
$ws.Range("E43").Value = "  +0.46%  "

# Row 44: update E44
$ws.Range("E44").Value = "  +0.81%  "

# Row 45: update D45, E45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0422"
$ws.Range("E45").Value = "  +1.39%  "

# Row 46: update D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.130"
$ws.Range("E46").Value = "  +1.25%  "

# Row 47: update E47
$ws.Range("E47").Value = "  +1.46%  "

# Row 48: update E48
$ws.Range("E48").Value = "  -0.15%  "

# Row 49: update E49
$ws.Range("E49").Value = "  +10.58%  "

# Row 50: update D50, E50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.85"
$ws.Range("E50").Value = "  -17.90%  "

# Row 51: update D51, E51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.14"
$ws.Range("E51").Value = "  +3.12%  "
